$d = $word.ActiveDocument

# Locate the paragraph that holds the "LOQ4031 ... (Requisito fraco)" text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOQ4031*Requisito fraco*") {
        $target = $p
    }
}

if ($target -ne $null) {
    # Immediately after the requirement paragraph there is an empty paragraph,
    # then the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph, then the
    # copyright/footer paragraph. Remove that trio (the footer content and the
    # extra blank line), leaving just the single blank paragraph that follows
    # them before the page break.
    $blankPara = $target.Next()
    $jupiterPara = $blankPara.Next()
    $copyrightPara = $jupiterPara.Next()

    $r = $d.Range($blankPara.Range.Start, $copyrightPara.Range.End)
    $r.Delete()
}
